$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the last two rows (old rows 6 and 7) - data is being cleaned up
$ws.Range("A6:A7").EntireRow.Delete()

# Update remaining data rows (2-5) with new values:
# Columns: A=Subscription ID, B=Resource Group, C=DB Name, D=SKU,
#          E=Replicas per Master, F=Shard Count, G=Shard Number,
#          H=Avg Ops/Sec, I=Used Memory (MB), J=Max Total Connections

$guid = "00000000-0000-0000-0000-000000000000"

# Assign string columns in SKU, Resource Group, DB Name, Subscription ID order
# (matching the order these new values were introduced in the source data),
# then fill in the remaining numeric columns.
$ws.Cells.Item(2,4).Value = "C0 Basic"
$ws.Cells.Item(3,4).Value = "P1 Premium"
$ws.Cells.Item(4,4).Value = "P1 Premium"
$ws.Cells.Item(5,4).Value = "P1 Premium"

$ws.Cells.Item(2,2).Value = "demo-rg"
$ws.Cells.Item(3,2).Value = "demo-rg"
$ws.Cells.Item(4,2).Value = "demo-rg"
$ws.Cells.Item(5,2).Value = "demo-rg"

$ws.Cells.Item(2,3).Value = "demo-1"
$ws.Cells.Item(3,3).Value = "demo-2"
$ws.Cells.Item(4,3).Value = "demo-2"
$ws.Cells.Item(5,3).Value = "demo-2"

$ws.Cells.Item(2,1).Value = $guid
$ws.Cells.Item(3,1).Value = $guid
$ws.Cells.Item(4,1).Value = $guid
$ws.Cells.Item(5,1).Value = $guid

$ws.Cells.Item(2,5).Value = 0
$ws.Cells.Item(2,6).Value = 1
$ws.Cells.Item(2,7).Value = 0
$ws.Cells.Item(2,8).Value = 5
$ws.Cells.Item(2,9).Value = 0.7
$ws.Cells.Item(2,10).Value = 19

$ws.Cells.Item(3,5).Value = 1
$ws.Cells.Item(3,6).Value = 3
$ws.Cells.Item(3,7).Value = 0
$ws.Cells.Item(3,8).Value = 36
$ws.Cells.Item(3,9).Value = 4568.99
$ws.Cells.Item(3,10).Value = 8

$ws.Cells.Item(4,5).Value = 1
$ws.Cells.Item(4,6).Value = 3
$ws.Cells.Item(4,7).Value = 1
$ws.Cells.Item(4,8).Value = 28
$ws.Cells.Item(4,9).Value = 4569.54
$ws.Cells.Item(4,10).Value = 7

$ws.Cells.Item(5,5).Value = 1
$ws.Cells.Item(5,6).Value = 3
$ws.Cells.Item(5,7).Value = 2
$ws.Cells.Item(5,8).Value = 44
$ws.Cells.Item(5,9).Value = 4571.06
$ws.Cells.Item(5,10).Value = 8

# Adjust column widths to reflect new (shorter) content, matching the
# best-fit widths Excel computed after the data changed.
$ws.Columns("B").ColumnWidth = 13.365885416666666
$ws.Columns("C").ColumnWidth = 7.631510416666667
$ws.Columns("D").ColumnWidth = 10.029947916666666

# Update the selection to match the final state
$ws.Range("A5").Select()
